$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: Anthony Edwards -> Collin Sexton / PG,SG / Utah Jazz
$ws.Range("A3").Value = "Collin Sexton"
$ws.Range("B3").Value = "PG,SG"
$ws.Range("C3").Value = "Utah Jazz"

# Row 6: Jayson Tatum -> Andrew Wiggins / SF,PF / Golden State Warriors
$ws.Range("A6").Value = "Andrew Wiggins"
$ws.Range("B6").Value = "SF,PF"
$ws.Range("C6").Value = "Golden State Warriors"

# Row 14: Collin Sexton -> D'Angelo Russell / PG / Los Angeles Lakers
$ws.Range("A14").Value = "D'Angelo Russell"
$ws.Range("B14").Value = "PG"
$ws.Range("C14").Value = "Los Angeles Lakers"

# Row 15: Andrew Wiggins -> Anthony Edwards / SG,SF / Minnesota Timberwolves
$ws.Range("A15").Value = "Anthony Edwards"
$ws.Range("B15").Value = "SG,SF"
$ws.Range("C15").Value = "Minnesota Timberwolves"

# Row 16: Draymond Green -> Jayson Tatum / SF,PF / Boston Celtics
$ws.Range("A16").Value = "Jayson Tatum"
$ws.Range("B16").Value = "SF,PF"
$ws.Range("C16").Value = "Boston Celtics"

# Row 18: Kyle Kuzma -> Draymond Green / PF,C / Golden State Warriors
$ws.Range("A18").Value = "Draymond Green"
$ws.Range("B18").Value = "PF,C"
$ws.Range("C18").Value = "Golden State Warriors"
